$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a "clean" number-looking string (e.g. "254.31").
# Excel auto-converts such text to a numeric cell on assignment, which would
# change both the stored type and the original "98.021.11"-style text layout
# (trailing zeros, etc.). Force Text format for the assignment, then restore
# the cell to the workbooks default (unstyled) look so no stray number format
# is left behind.

# Row 2
$ws.Range("D2").Value = '98.025.73'
$ws.Range("E2").Value = '  -0.48%  '

# Row 3
$ws.Range("D3").Value = '3.403.90'
$ws.Range("E3").Value = '  +1.22%  '

# Row 4
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '254.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.51%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '678.19'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +2.01%  '

# Row 7
$ws.Range("E7").Value = '  -6.05%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.432'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -7.20%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.05'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.04%  '

# Row 10
$ws.Range("E10").Value = '  -0.04%  '

# Row 11
$ws.Range("D11").Value = '3.402.16'
$ws.Range("E11").Value = '  +1.29%  '

# Row 12
$ws.Range("E12").Value = '  +0.61%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.07%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.28'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +11.25%  '

# Row 15
$ws.Range("D15").Value = '97.821.49'
$ws.Range("E15").Value = '  -1.64%  '

# Row 16
$ws.Range("E16").Value = '  -2.53%  '

# Row 17
$ws.Range("D17").Value = '4.037.94'
$ws.Range("E17").Value = '  +1.13%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.85'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +15.31%  '

# Row 19
$ws.Range("D19").Value = '3.419.49'
$ws.Range("E19").Value = '  +1.68%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.572'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +31.08%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.40'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.04%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.06'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.28%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.44'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -4.43%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '507.57'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.71%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000204'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.95%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.54'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.34%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '99.35'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.99%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.60'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.14%  '

# Row 29
$ws.Range("D29").Value = '3.582.35'
$ws.Range("E29").Value = '  +1.05%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.151'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.94%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.54'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.75%  '

# Row 32
$ws.Range("E32").Value = '  -0.04%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.195'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +2.07%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.60'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +22.26%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.34%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.570'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.26%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '29.41'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.33%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.53'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +13.61%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.96'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.19%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '536.30'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.76%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.153'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.69%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '24.71'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.06%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.870'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.24%  '

# Row 45
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0434'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.44%  '

# Row 46
$ws.Range("B46").Value = 'Cosmos'
$ws.Range("C46").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.98'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +12.39%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.76'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.53%  '

# Row 48
$ws.Range("E48").Value = '  +13.41%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.75'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +11.85%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '55.46'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.33%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '3.21'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -6.41%  '

